$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# =========================================================================
# STEP A: copy cell formatting (styles) into the new/changed cells, using
# pristine source cells so results don't depend on later value edits.
# =========================================================================

# Header row: G1/J1 take the "numeric header" style (same as F1); H1/I1/K1/L1
# take the "text header" style (same as A1).
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("A1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("L1").PasteSpecial(-4122)

# Data rows 2-8: G (absolute MPSP diff) gets a new style derived from H's
# (percentage) style but with a plain 2-decimal number format. I/K/L (the
# relative-diff / consistency columns) reuse H's existing style untouched.
# J (absolute GWP diff) reuses G's newly derived style.
for ($r = 2; $r -le 8; $r++) {
    $ws.Range("H$r").Copy()
    $ws.Range("I$r").PasteSpecial(-4122)
    $ws.Range("K$r").PasteSpecial(-4122)
    $ws.Range("L$r").PasteSpecial(-4122)

    $ws.Range("G$r").PasteSpecial(-4122)
    $ws.Range("G$r").NumberFormat = "0.00"

    $ws.Range("G$r").Copy()
    $ws.Range("J$r").PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# =========================================================================
# STEP B: values / formulas
# =========================================================================

# Header labels (order chosen to reproduce the shared-string table order)
$ws.Range("D1").Value = "GWP_original"
$ws.Range("E1").Value = "MPSP_wwt"
$ws.Range("F1").Value = "GWP_wwt"
$ws.Range("H1").Value = "MPSP_diff_rel"
$ws.Range("G1").Value = "MPSP_diff_abs"
$ws.Range("J1").Value = "GWP_diff_abs"
$ws.Range("K1").Value = "GWP_diff_rel"
$ws.Range("I1").Value = "consistency"
$ws.Range("L1").Value = "consistency"

# Data formulas
for ($r = 2; $r -le 8; $r++) {
    $ws.Range("G$r").Formula = "=(E$r-C$r)"
    $ws.Range("H$r").Formula = "=(E$r-C$r)/C$r"
    $ws.Range("I$r").Formula = "=IF(AND(ABS(G$r)>0.1, ABS(H$r)>0.1), FALSE, TRUE)"
    $ws.Range("J$r").Formula = "=IF(D$r<>`"NA`",(F$r-D$r), 0)"
    $ws.Range("K$r").Formula = "=IF(D$r<>`"NA`",(F$r-D$r)/D$r, 0)"
    $ws.Range("L$r").Formula = "=IF(AND(ABS(J$r)>0.1, ABS(K$r)>0.1), FALSE, TRUE)"
}

# =========================================================================
# STEP C: column widths (EntireColumn.ColumnWidth snaps to the nearest
# 1/6 of a character before Excel adds its ~0.833 char padding, so we pick
# the closest reachable value to each target width).
# =========================================================================
$ws.Range("G1").EntireColumn.ColumnWidth = 12.666666666666666
$ws.Range("H1").EntireColumn.ColumnWidth = 12
$ws.Range("I1").EntireColumn.ColumnWidth = 11.5
$ws.Range("J1").EntireColumn.ColumnWidth = 12.666666666666666
$ws.Range("K1").EntireColumn.ColumnWidth = 11.5
$ws.Range("L1").EntireColumn.ColumnWidth = 11.5

# =========================================================================
# STEP D: conditional formatting
# =========================================================================
$ws.Range("G2:H8").FormatConditions.Delete()

# Consistency (TRUE/FALSE) highlighting on the boolean check columns.
$fcTrue = $ws.Range("H2:I8").FormatConditions.Add(1, 3, "TRUE")
$fcTrue.Interior.Color = 13561798
$fcTrue.Font.Color = 24832

$fcFalse = $ws.Range("H2:I8").FormatConditions.Add(1, 3, "FALSE")
$fcFalse.Interior.Color = 13551615
$fcFalse.Font.Color = 393372

$fcTrue2 = $ws.Range("L2:L8").FormatConditions.Add(1, 3, "TRUE")
$fcTrue2.Interior.Color = 13561798
$fcTrue2.Font.Color = 24832

$fcFalse2 = $ws.Range("L2:L8").FormatConditions.Add(1, 3, "FALSE")
$fcFalse2.Interior.Color = 13551615
$fcFalse2.Font.Color = 393372

# Out-of-tolerance highlighting on the numeric diff columns.
$fcG = $ws.Range("G2:G8").FormatConditions.Add(1, 2, "-0.1", "0.1")
$fcG.Interior.Color = 13551615
$fcG.Font.Color = 393372

$fcJ = $ws.Range("J2:J8").FormatConditions.Add(1, 2, "-0.1", "0.1")
$fcJ.Interior.Color = 13551615
$fcJ.Font.Color = 393372

$fcK = $ws.Range("K2:K8").FormatConditions.Add(1, 2, "-0.1", "0.1")
$fcK.Interior.Color = 13551615
$fcK.Font.Color = 393372

# =========================================================================
# STEP E: selection / active cell
# =========================================================================
[void]$ws.Range("L1:L1048576").Select()
